$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A14").Value = "AddressBook"
$ws.Range("E14").Value = "qa"
$ws.Range("F14").Value = "Testing"
$ws.Range("I14").Value = "6 Walnut Valley Dr"
$ws.Range("J14").Value = "Little Rock"
$ws.Range("K14").Value = "Arkansas"
$ws.Range("L14").Value = "72211"
$ws.Range("M14").Value = "9999999999"
